# Clear the "No. of Sites/bldg ..." breakdown columns (AB:AK) and the
# "DIFFERENCE" column (AM) for the data rows, keeping only the
# "PREVIOUS ACCOMPLISHMENT" column (AL) as the up-to-date status value,
# per the May accomplishment report update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 15; $row++) {
    $ws.Range("AB${row}:AK${row}").ClearContents()
    $ws.Range("AM${row}").ClearContents()
}
